# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect refreshed data from the gh-pages generator.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 102
    $ws.Range("F3").Value = 923
    $ws.Range("F4").Value = 217
}
